$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 2214.2856
$ws.Cells.Item(40, 9).Value = 2000
$ws.Cells.Item(40, 11).Value = 2000
$ws.Cells.Item(40, 13).Value = -1825
$ws.Cells.Item(62, 8).Value = 7025.625
$ws.Cells.Item(62, 9).Value = 5401.6665
$ws.Cells.Item(62, 10).Value = 8000
$ws.Cells.Item(62, 11).Value = 5401.6665
$ws.Cells.Item(62, 12).Value = 8000
$ws.Cells.Item(62, 13).Value = -4777.6665
$ws.Cells.Item(62, 14).Value = -9248
$ws.Cells.Item(65, 8).Value = 7025.625
$ws.Cells.Item(65, 9).Value = 5401.6665
$ws.Cells.Item(65, 10).Value = 8000
$ws.Cells.Item(65, 11).Value = 27008.3325
$ws.Cells.Item(65, 12).Value = 40000
$ws.Cells.Item(65, 13).Value = -23888.3325
$ws.Cells.Item(65, 14).Value = -46240
$ws.Cells.Item(92, 8).Value = 673.2778
$ws.Cells.Item(92, 9).Value = 554.3333
$ws.Cells.Item(92, 10).Value = 1268
$ws.Cells.Item(92, 11).Value = 554.3333
$ws.Cells.Item(92, 12).Value = 1268
$ws.Cells.Item(92, 13).Value = 693.6667
$ws.Cells.Item(92, 14).Value = -3764
$ws.Cells.Item(112, 8).Value = 2411
$ws.Cells.Item(112, 10).Value = 2452.1
$ws.Cells.Item(112, 12).Value = 7356.299999999999
$ws.Cells.Item(112, 14).Value = -9572.299999999999
$ws.Cells.Item(113, 8).Value = 2179.8
$ws.Cells.Item(113, 9).Value = 2366.6667
$ws.Cells.Item(113, 10).Value = 1899.5
$ws.Cells.Item(113, 11).Value = 2366.6667
$ws.Cells.Item(113, 12).Value = 1899.5
$ws.Cells.Item(113, 13).Value = 887.3332999999998
$ws.Cells.Item(113, 14).Value = -8407.5
$ws.Cells.Item(116, 8).Value = 4544.7
$ws.Cells.Item(116, 9).Value = 4528.143
$ws.Cells.Item(116, 10).Value = 4583.3335
$ws.Cells.Item(116, 11).Value = 4528.143
$ws.Cells.Item(116, 12).Value = 4583.3335
$ws.Cells.Item(116, 13).Value = -1086.143
$ws.Cells.Item(116, 14).Value = -11467.3335
$ws.Cells.Item(132, 8).Value = 2628.1667
$ws.Cells.Item(132, 9).Value = 2594.3635
$ws.Cells.Item(132, 11).Value = 7783.0905
$ws.Cells.Item(132, 13).Value = -5253.0905
$ws.Cells.Item(138, 8).Value = 4179.0454
$ws.Cells.Item(138, 10).Value = 4766.8057
$ws.Cells.Item(138, 12).Value = 14300.4171
$ws.Cells.Item(138, 14).Value = -24580.4171

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 11
$ws.Cells.Item(2, 9).Value = 11
$ws.Cells.Item(2, 11).Value = 11
$ws.Cells.Item(2, 13).Value = 102
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(61, 8).Value = 4560.6
$ws.Cells.Item(61, 9).Value = 4536.6665
$ws.Cells.Item(61, 11).Value = 4536.6665
$ws.Cells.Item(61, 13).Value = -4324.6665
$ws.Cells.Item(116, 8).Value = 11
$ws.Cells.Item(116, 9).Value = 11
$ws.Cells.Item(116, 11).Value = 11
$ws.Cells.Item(116, 13).Value = 2283
$ws.Cells.Item(136, 8).Value = 4560.6
$ws.Cells.Item(136, 9).Value = 4536.6665
$ws.Cells.Item(136, 11).Value = 13609.9995
$ws.Cells.Item(136, 13).Value = -11059.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 11
$ws.Cells.Item(3, 9).Value = 11
$ws.Cells.Item(3, 11).Value = 11
$ws.Cells.Item(3, 13).Value = 103
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(134, 8).Value = 3092.3333
$ws.Cells.Item(134, 9).Value = 2845.5557
$ws.Cells.Item(134, 11).Value = 8536.667099999999
$ws.Cells.Item(134, 13).Value = -6001.667099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 184.22223
$ws.Cells.Item(7, 9).Value = 165.64706
$ws.Cells.Item(7, 11).Value = 165.64706
$ws.Cells.Item(7, 13).Value = -52.64706000000001
$ws.Cells.Item(10, 8).Value = 291
$ws.Cells.Item(10, 9).Value = 339.2
$ws.Cells.Item(10, 11).Value = 339.2
$ws.Cells.Item(10, 13).Value = -200.2
$ws.Cells.Item(58, 8).Value = 3603.3914
$ws.Cells.Item(58, 9).Value = 1753.4166
$ws.Cells.Item(58, 11).Value = 1753.4166
$ws.Cells.Item(58, 13).Value = -1550.4166
$ws.Cells.Item(99, 8).Value = 11092.946
$ws.Cells.Item(99, 10).Value = 13407.5
$ws.Cells.Item(99, 12).Value = 13407.5
$ws.Cells.Item(99, 14).Value = -16403.5
$ws.Cells.Item(126, 8).Value = 11092.946
$ws.Cells.Item(126, 10).Value = 13407.5
$ws.Cells.Item(126, 12).Value = 40222.5
$ws.Cells.Item(126, 14).Value = -45162.5
$ws.Cells.Item(136, 8).Value = 3603.3914
$ws.Cells.Item(136, 9).Value = 1753.4166
$ws.Cells.Item(136, 11).Value = 5260.2498
$ws.Cells.Item(136, 13).Value = -2710.2498

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 47665.19
$ws.Cells.Item(2, 9).Value = 76948.16
$ws.Cells.Item(2, 10).Value = 80.375
$ws.Cells.Item(2, 11).Value = 461688.96
$ws.Cells.Item(2, 12).Value = 482.25
$ws.Cells.Item(2, 13).Value = -461575.96
$ws.Cells.Item(2, 14).Value = -708.25
$ws.Cells.Item(9, 8).Value = 750.1667
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 221
$ws.Cells.Item(17, 8).Value = 71
$ws.Cells.Item(17, 9).Value = 68.25
$ws.Cells.Item(17, 10).Value = 74.666664
$ws.Cells.Item(17, 11).Value = 204.75
$ws.Cells.Item(17, 12).Value = 223.999992
$ws.Cells.Item(17, 13).Value = -35.75
$ws.Cells.Item(17, 14).Value = -561.999992
$ws.Cells.Item(19, 8).Value = 19.333334
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(24, 8).Value = 98
$ws.Cells.Item(24, 9).Value = 98
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 294
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = -64
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(35, 8).Value = 85
$ws.Cells.Item(35, 9).Value = 20
$ws.Cells.Item(35, 10).Value = 117.5
$ws.Cells.Item(35, 11).Value = 60
$ws.Cells.Item(35, 12).Value = 352.5
$ws.Cells.Item(35, 13).Value = 228
$ws.Cells.Item(35, 14).Value = -928.5
$ws.Cells.Item(107, 8).Value = 490.9355
$ws.Cells.Item(107, 10).Value = 482.60715
$ws.Cells.Item(107, 12).Value = 1447.82145
$ws.Cells.Item(107, 14).Value = -5287.821449999999
$ws.Cells.Item(119, 8).Value = 2831
$ws.Cells.Item(119, 9).Value = 2831
$ws.Cells.Item(119, 11).Value = 8493
$ws.Cells.Item(119, 13).Value = -3655
$ws.Cells.Item(129, 8).Value = 2421.2
$ws.Cells.Item(129, 9).Value = 996.6667
$ws.Cells.Item(129, 10).Value = 3031.7144
$ws.Cells.Item(129, 11).Value = 2990.0001
$ws.Cells.Item(129, 12).Value = 9095.143199999999
$ws.Cells.Item(129, 13).Value = 2009.9999
$ws.Cells.Item(129, 14).Value = -19095.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 992.25
$ws.Cells.Item(97, 9).Value = 904.5
$ws.Cells.Item(97, 10).Value = 1080
$ws.Cells.Item(97, 11).Value = 904.5
$ws.Cells.Item(97, 12).Value = 1080
$ws.Cells.Item(97, 13).Value = -408.5
$ws.Cells.Item(97, 14).Value = -2072
$ws.Cells.Item(102, 8).Value = 310.8
$ws.Cells.Item(102, 9).Value = 310.8
$ws.Cells.Item(102, 11).Value = 310.8
$ws.Cells.Item(102, 13).Value = 1311.2
$ws.Cells.Item(130, 8).Value = 100000
$ws.Cells.Item(130, 10).Value = 100000
$ws.Cells.Item(130, 12).Value = 100000
$ws.Cells.Item(130, 14).Value = -110040

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2186.7144
$ws.Cells.Item(40, 9).Value = 2186.7144
$ws.Cells.Item(40, 11).Value = 2186.7144
$ws.Cells.Item(40, 13).Value = -2050.7144
$ws.Cells.Item(46, 9).Value = 2499.75
$ws.Cells.Item(46, 10).Value = 3143.6667
$ws.Cells.Item(46, 11).Value = 2499.75
$ws.Cells.Item(46, 12).Value = 3143.6667
$ws.Cells.Item(46, 13).Value = -2311.75
$ws.Cells.Item(46, 14).Value = -3519.6667
$ws.Cells.Item(99, 8).Value = 21925.334
$ws.Cells.Item(99, 9).Value = 21925.334
$ws.Cells.Item(99, 11).Value = 21925.334
$ws.Cells.Item(99, 13).Value = -18930.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(74, 8).Value = 18838
$ws.Cells.Item(74, 9).Value = 4001
$ws.Cells.Item(74, 10).Value = 22547.25
$ws.Cells.Item(74, 11).Value = 4001
$ws.Cells.Item(74, 12).Value = 22547.25
$ws.Cells.Item(74, 13).Value = -3065
$ws.Cells.Item(74, 14).Value = -24419.25
$ws.Cells.Item(77, 8).Value = 18838
$ws.Cells.Item(77, 9).Value = 4001
$ws.Cells.Item(77, 10).Value = 22547.25
$ws.Cells.Item(77, 11).Value = 12003
$ws.Cells.Item(77, 12).Value = 67641.75
$ws.Cells.Item(77, 13).Value = -7323
$ws.Cells.Item(77, 14).Value = -77001.75
$ws.Cells.Item(136, 8).Value = 1857.2858
$ws.Cells.Item(136, 9).Value = 1615.5769
$ws.Cells.Item(136, 10).Value = 4999.5
$ws.Cells.Item(136, 11).Value = 4846.7307
$ws.Cells.Item(136, 12).Value = 14998.5
$ws.Cells.Item(136, 13).Value = -2296.7307
$ws.Cells.Item(136, 14).Value = -20098.5
